$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F (shifts MaxLevel/Value(1)/LevelUpAmount right to G/H/I)
$ws.Columns("F:F").Insert()

# New header + values for the inserted "StartLevel" column
$ws.Range("F1").Value = "StartLevel"
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1

# Match the column width used for the new column in the target workbook
$ws.Columns("F:F").ColumnWidth = 20.69921875

# Update sheet view (scrolled right a bit, new selection)
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("F15").Select()
